# Update "想去人数" (F column) figures across the sheets to the
# freshly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet "展览")
$wsExhibit.Range("F3").Value  = 2612
$wsExhibit.Range("F5").Value  = 917
$wsExhibit.Range("F7").Value  = 1923
$wsExhibit.Range("F8").Value  = 1763
$wsExhibit.Range("F9").Value  = 200
$wsExhibit.Range("F11").Value = 2423
$wsExhibit.Range("F12").Value = 522
$wsExhibit.Range("F13").Value = 197
$wsExhibit.Range("F16").Value = 113
$wsExhibit.Range("F17").Value = 102
$wsExhibit.Range("F18").Value = 8894
$wsExhibit.Range("F19").Value = 51
$wsExhibit.Range("F20").Value = 6920
$wsExhibit.Range("F21").Value = 11279
$wsExhibit.Range("F25").Value = 311
$wsExhibit.Range("F26").Value = 537
$wsExhibit.Range("F27").Value = 2468
$wsExhibit.Range("F28").Value = 213
$wsExhibit.Range("F29").Value = 185
$wsExhibit.Range("F30").Value = 2362
$wsExhibit.Range("F31").Value = 555
$wsExhibit.Range("F32").Value = 34
$wsExhibit.Range("F33").Value = 4477
$wsExhibit.Range("F34").Value = 709
$wsExhibit.Range("F35").Value = 315
$wsExhibit.Range("F36").Value = 26
$wsExhibit.Range("F37").Value = 481

# 演出 (sheet "演出")
$wsShow.Range("F8").Value  = 1180
$wsShow.Range("F23").Value = 5

# 本地生活 (sheet "本地生活")
$wsLocal.Range("F3").Value = 619
$wsLocal.Range("F5").Value = 121

# 全部类型 (sheet "全部类型") - aggregated view of all the above
$wsAll.Range("F3").Value  = 619
$wsAll.Range("F6").Value  = 2612
$wsAll.Range("F8").Value  = 917
$wsAll.Range("F10").Value = 1923
$wsAll.Range("F12").Value = 1763
$wsAll.Range("F14").Value = 200
$wsAll.Range("F15").Value = 2423
$wsAll.Range("F17").Value = 522
$wsAll.Range("F18").Value = 197
$wsAll.Range("F21").Value = 113
$wsAll.Range("F22").Value = 102
$wsAll.Range("F23").Value = 8894
$wsAll.Range("F24").Value = 51
$wsAll.Range("F25").Value = 6920
$wsAll.Range("F26").Value = 11279
$wsAll.Range("F32").Value = 537
$wsAll.Range("F36").Value = 213
$wsAll.Range("F37").Value = 185
$wsAll.Range("F38").Value = 34
$wsAll.Range("F39").Value = 4477
$wsAll.Range("F46").Value = 481
$wsAll.Range("F49").Value = 5
